# Update the LR-pairs sheet with newly recomputed TPM-based values.
# The original data only contained "FAPs" as the sending cluster (rows 2-7).
# The refreshed export adds an "ECs" sending-cluster block (now rows 2-7)
# ahead of the original "FAPs" block, which moves down to rows 8-13 and
# gets refreshed numeric values as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

  # Row 2
  $ws.Cells.Item(2, 1).Value = 'ECs'
  $ws.Cells.Item(2, 2).Value = 'Wnt5a'
  $ws.Cells.Item(2, 3).Value = 'Mcam'
  $ws.Cells.Item(2, 4).Value = 'ECs'
  $ws.Cells.Item(2, 5).Value = 2
  $ws.Cells.Item(2, 6).Value = 1
  $ws.Cells.Item(2, 7).Value = 0.125615
  $ws.Cells.Item(2, 8).Value = 0.25123
  $ws.Cells.Item(2, 9).Value = 0.02647478672532295
  $ws.Cells.Item(2, 10).Value = 0.01780700335556722
  $ws.Cells.Item(2, 11).Value = 2
  $ws.Cells.Item(2, 12).Value = 1
  $ws.Cells.Item(2, 13).Value = 138.2205735
  $ws.Cells.Item(2, 14).Value = 276.441147
  $ws.Cells.Item(2, 15).Value = 0.626961773180729
  $ws.Cells.Item(2, 16).Value = 0.6221415178520601
  $ws.Cells.Item(2, 17).Value = 17.3625773402025
  $ws.Cells.Item(2, 18).Value = 69.45030936081
  $ws.Cells.Item(2, 19).Value = 0.01659867922989011
  $ws.Cells.Item(2, 20).Value = 0.01107847609602932

  # Row 3
  $ws.Cells.Item(3, 1).Value = 'ECs'
  $ws.Cells.Item(3, 2).Value = 'Wnt5a'
  $ws.Cells.Item(3, 3).Value = 'Mcam'
  $ws.Cells.Item(3, 4).Value = 'FAPs'
  $ws.Cells.Item(3, 5).Value = 2
  $ws.Cells.Item(3, 6).Value = 1
  $ws.Cells.Item(3, 7).Value = 0.125615
  $ws.Cells.Item(3, 8).Value = 0.25123
  $ws.Cells.Item(3, 9).Value = 0.02647478672532295
  $ws.Cells.Item(3, 10).Value = 0.01780700335556722
  $ws.Cells.Item(3, 11).Value = 3
  $ws.Cells.Item(3, 12).Value = 1
  $ws.Cells.Item(3, 13).Value = 0.9119573333333332
  $ws.Cells.Item(3, 14).Value = 2.735872
  $ws.Cells.Item(3, 15).Value = 0.004136593940350246
  $ws.Cells.Item(3, 16).Value = 0.006157185994923365
  $ws.Cells.Item(3, 17).Value = 0.1145555204266667
  $ws.Cells.Item(3, 18).Value = 0.6873331225599999
  $ws.Cells.Item(3, 19).Value = 0.0001095154423400361
  $ws.Cells.Item(3, 20).Value = 0.0001096410316724519

  # Row 4
  $ws.Cells.Item(4, 1).Value = 'ECs'
  $ws.Cells.Item(4, 2).Value = 'Wnt5a'
  $ws.Cells.Item(4, 3).Value = 'Mcam'
  $ws.Cells.Item(4, 4).Value = 'Inflammatory-Mac'
  $ws.Cells.Item(4, 5).Value = 2
  $ws.Cells.Item(4, 6).Value = 1
  $ws.Cells.Item(4, 7).Value = 0.125615
  $ws.Cells.Item(4, 8).Value = 0.25123
  $ws.Cells.Item(4, 9).Value = 0.02647478672532295
  $ws.Cells.Item(4, 10).Value = 0.01780700335556722
  $ws.Cells.Item(4, 11).Value = 2
  $ws.Cells.Item(4, 12).Value = 0.6666666666666666
  $ws.Cells.Item(4, 13).Value = 0.3000813333333334
  $ws.Cells.Item(4, 14).Value = 0.900244
  $ws.Cells.Item(4, 15).Value = 0.001361154277406497
  $ws.Cells.Item(4, 16).Value = 0.002026034020894907
  $ws.Cells.Item(4, 17).Value = 0.03769471668666667
  $ws.Cells.Item(4, 18).Value = 0.22616830012
  $ws.Cells.Item(4, 19).Value = 0.00003603626919459808
  $ws.Cells.Item(4, 20).Value = 0.00003607759460856895

  # Row 5
  $ws.Cells.Item(5, 1).Value = 'ECs'
  $ws.Cells.Item(5, 2).Value = 'Wnt5a'
  $ws.Cells.Item(5, 3).Value = 'Mcam'
  $ws.Cells.Item(5, 4).Value = 'MuSCs'
  $ws.Cells.Item(5, 5).Value = 2
  $ws.Cells.Item(5, 6).Value = 1
  $ws.Cells.Item(5, 7).Value = 0.125615
  $ws.Cells.Item(5, 8).Value = 0.25123
  $ws.Cells.Item(5, 9).Value = 0.02647478672532295
  $ws.Cells.Item(5, 10).Value = 0.01780700335556722
  $ws.Cells.Item(5, 11).Value = 2
  $ws.Cells.Item(5, 12).Value = 1
  $ws.Cells.Item(5, 13).Value = 78.82415950000001
  $ws.Cells.Item(5, 14).Value = 157.648319
  $ws.Cells.Item(5, 15).Value = 0.357542539132936
  $ws.Cells.Item(5, 16).Value = 0.3547936533105391
  $ws.Cells.Item(5, 17).Value = 9.901496795592502
  $ws.Cells.Item(5, 18).Value = 39.60598718237001
  $ws.Cells.Item(5, 19).Value = 0.009465862468774918
  $ws.Cells.Item(5, 20).Value = 0.006317811775034723

  # Row 6
  $ws.Cells.Item(6, 1).Value = 'ECs'
  $ws.Cells.Item(6, 2).Value = 'Wnt5a'
  $ws.Cells.Item(6, 3).Value = 'Mcam'
  $ws.Cells.Item(6, 4).Value = 'Neutrophils'
  $ws.Cells.Item(6, 5).Value = 2
  $ws.Cells.Item(6, 6).Value = 1
  $ws.Cells.Item(6, 7).Value = 0.125615
  $ws.Cells.Item(6, 8).Value = 0.25123
  $ws.Cells.Item(6, 9).Value = 0.02647478672532295
  $ws.Cells.Item(6, 10).Value = 0.01780700335556722
  $ws.Cells.Item(6, 11).Value = 3
  $ws.Cells.Item(6, 12).Value = 1
  $ws.Cells.Item(6, 13).Value = 1.659958666666667
  $ws.Cells.Item(6, 14).Value = 4.979876
  $ws.Cells.Item(6, 15).Value = 0.007529491469372698
  $ws.Cells.Item(6, 16).Value = 0.01120740398807217
  $ws.Cells.Item(6, 17).Value = 0.2085157079133333
  $ws.Cells.Item(6, 18).Value = 1.25109424748
  $ws.Cells.Item(6, 19).Value = 0.0001993416808017807
  $ws.Cells.Item(6, 20).Value = 0.0001995702804227986

  # Row 7
  $ws.Cells.Item(7, 1).Value = 'ECs'
  $ws.Cells.Item(7, 2).Value = 'Wnt5a'
  $ws.Cells.Item(7, 3).Value = 'Mcam'
  $ws.Cells.Item(7, 4).Value = 'Resolving-Mac'
  $ws.Cells.Item(7, 5).Value = 2
  $ws.Cells.Item(7, 6).Value = 1
  $ws.Cells.Item(7, 7).Value = 0.125615
  $ws.Cells.Item(7, 8).Value = 0.25123
  $ws.Cells.Item(7, 9).Value = 0.02647478672532295
  $ws.Cells.Item(7, 10).Value = 0.01780700335556722
  $ws.Cells.Item(7, 11).Value = 3
  $ws.Cells.Item(7, 12).Value = 1
  $ws.Cells.Item(7, 13).Value = 0.5441963333333334
  $ws.Cells.Item(7, 14).Value = 1.632589
  $ws.Cells.Item(7, 15).Value = 0.002468447999205544
  $ws.Cells.Item(7, 16).Value = 0.003674204833510465
  $ws.Cells.Item(7, 17).Value = 0.06835922241166667
  $ws.Cells.Item(7, 18).Value = 0.41015533447
  $ws.Cells.Item(7, 19).Value = 0.00006535163432151693
  $ws.Cells.Item(7, 20).Value = 0.00006542657779936216

  # Row 8
  $ws.Cells.Item(8, 1).Value = 'FAPs'
  $ws.Cells.Item(8, 2).Value = 'Wnt5a'
  $ws.Cells.Item(8, 3).Value = 'Mcam'
  $ws.Cells.Item(8, 4).Value = 'ECs'
  $ws.Cells.Item(8, 5).Value = 3
  $ws.Cells.Item(8, 6).Value = 1
  $ws.Cells.Item(8, 7).Value = 4.619088000000001
  $ws.Cells.Item(8, 8).Value = 13.857264
  $ws.Cells.Item(8, 9).Value = 0.9735252132746771
  $ws.Cells.Item(8, 10).Value = 0.9821929966444328
  $ws.Cells.Item(8, 11).Value = 2
  $ws.Cells.Item(8, 12).Value = 1
  $ws.Cells.Item(8, 13).Value = 138.2205735
  $ws.Cells.Item(8, 14).Value = 276.441147
  $ws.Cells.Item(8, 15).Value = 0.626961773180729
  $ws.Cells.Item(8, 16).Value = 0.6221415178520601
  $ws.Cells.Item(8, 17).Value = 638.4529924069681
  $ws.Cells.Item(8, 18).Value = 3830.717954441808
  $ws.Cells.Item(8, 19).Value = 0.6103630939508389
  $ws.Cells.Item(8, 20).Value = 0.6110630417560308

  # Row 9
  $ws.Cells.Item(9, 1).Value = 'FAPs'
  $ws.Cells.Item(9, 2).Value = 'Wnt5a'
  $ws.Cells.Item(9, 3).Value = 'Mcam'
  $ws.Cells.Item(9, 4).Value = 'FAPs'
  $ws.Cells.Item(9, 5).Value = 3
  $ws.Cells.Item(9, 6).Value = 1
  $ws.Cells.Item(9, 7).Value = 4.619088000000001
  $ws.Cells.Item(9, 8).Value = 13.857264
  $ws.Cells.Item(9, 9).Value = 0.9735252132746771
  $ws.Cells.Item(9, 10).Value = 0.9821929966444328
  $ws.Cells.Item(9, 11).Value = 3
  $ws.Cells.Item(9, 12).Value = 1
  $ws.Cells.Item(9, 13).Value = 0.9119573333333332
  $ws.Cells.Item(9, 14).Value = 2.735872
  $ws.Cells.Item(9, 15).Value = 0.004136593940350246
  $ws.Cells.Item(9, 16).Value = 0.006157185994923365
  $ws.Cells.Item(9, 17).Value = 4.212411174912
  $ws.Cells.Item(9, 18).Value = 37.911700574208
  $ws.Cells.Item(9, 19).Value = 0.00402707849801021
  $ws.Cells.Item(9, 20).Value = 0.006047544963250913

  # Row 10
  $ws.Cells.Item(10, 1).Value = 'FAPs'
  $ws.Cells.Item(10, 2).Value = 'Wnt5a'
  $ws.Cells.Item(10, 3).Value = 'Mcam'
  $ws.Cells.Item(10, 4).Value = 'Inflammatory-Mac'
  $ws.Cells.Item(10, 5).Value = 3
  $ws.Cells.Item(10, 6).Value = 1
  $ws.Cells.Item(10, 7).Value = 4.619088000000001
  $ws.Cells.Item(10, 8).Value = 13.857264
  $ws.Cells.Item(10, 9).Value = 0.9735252132746771
  $ws.Cells.Item(10, 10).Value = 0.9821929966444328
  $ws.Cells.Item(10, 11).Value = 2
  $ws.Cells.Item(10, 12).Value = 0.6666666666666666
  $ws.Cells.Item(10, 13).Value = 0.3000813333333334
  $ws.Cells.Item(10, 14).Value = 0.900244
  $ws.Cells.Item(10, 15).Value = 0.001361154277406497
  $ws.Cells.Item(10, 16).Value = 0.002026034020894907
  $ws.Cells.Item(10, 17).Value = 1.386102085824
  $ws.Cells.Item(10, 18).Value = 12.474918772416
  $ws.Cells.Item(10, 19).Value = 0.001325118008211899
  $ws.Cells.Item(10, 20).Value = 0.001989956426286338

  # Row 11
  $ws.Cells.Item(11, 1).Value = 'FAPs'
  $ws.Cells.Item(11, 2).Value = 'Wnt5a'
  $ws.Cells.Item(11, 3).Value = 'Mcam'
  $ws.Cells.Item(11, 4).Value = 'MuSCs'
  $ws.Cells.Item(11, 5).Value = 3
  $ws.Cells.Item(11, 6).Value = 1
  $ws.Cells.Item(11, 7).Value = 4.619088000000001
  $ws.Cells.Item(11, 8).Value = 13.857264
  $ws.Cells.Item(11, 9).Value = 0.9735252132746771
  $ws.Cells.Item(11, 10).Value = 0.9821929966444328
  $ws.Cells.Item(11, 11).Value = 2
  $ws.Cells.Item(11, 12).Value = 1
  $ws.Cells.Item(11, 13).Value = 78.82415950000001
  $ws.Cells.Item(11, 14).Value = 157.648319
  $ws.Cells.Item(11, 15).Value = 0.357542539132936
  $ws.Cells.Item(11, 16).Value = 0.3547936533105391
  $ws.Cells.Item(11, 17).Value = 364.0957292565361
  $ws.Cells.Item(11, 18).Value = 2184.574375539216
  $ws.Cells.Item(11, 19).Value = 0.3480766766641611
  $ws.Cells.Item(11, 20).Value = 0.3484758415355044

  # Row 12
  $ws.Cells.Item(12, 1).Value = 'FAPs'
  $ws.Cells.Item(12, 2).Value = 'Wnt5a'
  $ws.Cells.Item(12, 3).Value = 'Mcam'
  $ws.Cells.Item(12, 4).Value = 'Neutrophils'
  $ws.Cells.Item(12, 5).Value = 3
  $ws.Cells.Item(12, 6).Value = 1
  $ws.Cells.Item(12, 7).Value = 4.619088000000001
  $ws.Cells.Item(12, 8).Value = 13.857264
  $ws.Cells.Item(12, 9).Value = 0.9735252132746771
  $ws.Cells.Item(12, 10).Value = 0.9821929966444328
  $ws.Cells.Item(12, 11).Value = 3
  $ws.Cells.Item(12, 12).Value = 1
  $ws.Cells.Item(12, 13).Value = 1.659958666666667
  $ws.Cells.Item(12, 14).Value = 4.979876
  $ws.Cells.Item(12, 15).Value = 0.007529491469372698
  $ws.Cells.Item(12, 16).Value = 0.01120740398807217
  $ws.Cells.Item(12, 17).Value = 7.667495157696001
  $ws.Cells.Item(12, 18).Value = 69.007456419264
  $ws.Cells.Item(12, 19).Value = 0.007330149788570918
  $ws.Cells.Item(12, 20).Value = 0.01100783370764937

  # Row 13
  $ws.Cells.Item(13, 1).Value = 'FAPs'
  $ws.Cells.Item(13, 2).Value = 'Wnt5a'
  $ws.Cells.Item(13, 3).Value = 'Mcam'
  $ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
  $ws.Cells.Item(13, 5).Value = 3
  $ws.Cells.Item(13, 6).Value = 1
  $ws.Cells.Item(13, 7).Value = 4.619088000000001
  $ws.Cells.Item(13, 8).Value = 13.857264
  $ws.Cells.Item(13, 9).Value = 0.9735252132746771
  $ws.Cells.Item(13, 10).Value = 0.9821929966444328
  $ws.Cells.Item(13, 11).Value = 3
  $ws.Cells.Item(13, 12).Value = 1
  $ws.Cells.Item(13, 13).Value = 0.5441963333333334
  $ws.Cells.Item(13, 14).Value = 1.632589
  $ws.Cells.Item(13, 15).Value = 0.002468447999205544
  $ws.Cells.Item(13, 16).Value = 0.003674204833510465
  $ws.Cells.Item(13, 17).Value = 2.513690752944
  $ws.Cells.Item(13, 18).Value = 22.623216776496
  $ws.Cells.Item(13, 19).Value = 0.002403096364884027
  $ws.Cells.Item(13, 20).Value = 0.003608778255711103
